$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Report")

# Header row
$ws.Range("A1").Value = "Student"
$ws.Range("B1").Value = "Grade"

# Clear any leftover data in column C (previously used)
$ws.Range("C1:C2").ClearContents()

# Student data rows
$ws.Range("A2").Value = "Tom"
$ws.Range("B2").Value = 4

$ws.Range("A3").Value = "Oscar"
$ws.Range("B3").Value = 4

$ws.Range("A4").Value = "Jay"
$ws.Range("B4").Value = 4

$ws.Range("A5").Value = "John"
$ws.Range("B5").Value = 4
